$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.00246799999999
$ws.Range("H2").Value = 70.00493599999999
$ws.Range("I2").Value = 0.2259773487303484
$ws.Range("J2").Value = 0.1743393500109074
$ws.Range("M2").Value = 0.2253945
$ws.Range("N2").Value = 0.450789
$ws.Range("O2").Value = 0.07045953565491116
$ws.Range("P2").Value = 0.05124470957648078
$ws.Range("Q2").Value = 7.889363773625998
$ws.Range("R2").Value = 31.55745509450399
$ws.Range("S2").Value = 0.01592225906006828
$ws.Range("T2").Value = 0.00893396935906138
$ws.Range("G3").Value = 35.00246799999999
$ws.Range("H3").Value = 70.00493599999999
$ws.Range("I3").Value = 0.2259773487303484
$ws.Range("J3").Value = 0.1743393500109074
$ws.Range("O3").Value = 0.7485640341142095
$ws.Range("P3").Value = 0.8166377943957063
$ws.Range("Q3").Value = 83.81681653288533
$ws.Range("R3").Value = 502.900899197312
$ws.Range("S3").Value = 0.1691585157840232
$ws.Range("T3").Value = 0.1423721022692884
$ws.Range("G4").Value = 35.00246799999999
$ws.Range("H4").Value = 70.00493599999999
$ws.Range("I4").Value = 0.2259773487303484
$ws.Range("J4").Value = 0.1743393500109074
$ws.Range("M4").Value = 0.574578
$ws.Range("N4").Value = 1.149156
$ws.Range("O4").Value = 0.1796161799756762
$ws.Range("P4").Value = 0.1306335458009631
$ws.Range("Q4").Value = 20.111648058504
$ws.Range("R4").Value = 80.446592234016
$ws.Range("S4").Value = 0.04058918813997641
$ws.Range("T4").Value = 0.02277456746456
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 35.00246799999999
$ws.Range("H5").Value = 70.00493599999999
$ws.Range("I5").Value = 0.2259773487303484
$ws.Range("J5").Value = 0.1743393500109074
$ws.Range("M5").Value = 0.004351333333333333
$ws.Range("N5").Value = 0.013054
$ws.Range("O5").Value = 0.001360250255203225
$ws.Range("P5").Value = 0.001483950226849768
$ws.Range("Q5").Value = 0.1523074057573333
$ws.Range("R5").Value = 0.9138444345439998
$ws.Range("S5").Value = 0.0003073857462806047
$ws.Range("T5").Value = 0.0002587109179975271
$ws.Range("G6").Value = 65.63887766666667
$ws.Range("I6").Value = 0.4237672483194333
$ws.Range("J6").Value = 0.4903985313772218
$ws.Range("M6").Value = 0.2253945
$ws.Range("N6").Value = 0.450789
$ws.Range("O6").Value = 0.07045953565491116
$ws.Range("P6").Value = 0.05124470957648078
$ws.Range("Q6").Value = 14.7946420122395
$ws.Range("R6").Value = 88.76785207343701
$ws.Range("S6").Value = 0.0298584435423467
$ws.Range("T6").Value = 0.02513033031715843
$ws.Range("G7").Value = 65.63887766666667
$ws.Range("I7").Value = 0.4237672483194333
$ws.Range("J7").Value = 0.4903985313772218
$ws.Range("O7").Value = 0.7485640341142095
$ws.Range("P7").Value = 0.8166377943957063
$ws.Range("S7").Value = 0.3172169209274729
$ws.Range("T7").Value = 0.400477975038788
$ws.Range("G8").Value = 65.63887766666667
$ws.Range("I8").Value = 0.4237672483194333
$ws.Range("J8").Value = 0.4903985313772218
$ws.Range("M8").Value = 0.574578
$ws.Range("N8").Value = 1.149156
$ws.Range("O8").Value = 0.1796161799756762
$ws.Range("P8").Value = 0.1306335458009631
$ws.Range("Q8").Value = 37.714655051958
$ws.Range("R8").Value = 226.287930311748
$ws.Range("S8").Value = 0.07611545434194039
$ws.Range("T8").Value = 0.06406249900939134
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 65.63887766666667
$ws.Range("I9").Value = 0.4237672483194333
$ws.Range("J9").Value = 0.4903985313772218
$ws.Range("M9").Value = 0.004351333333333333
$ws.Range("N9").Value = 0.013054
$ws.Range("O9").Value = 0.001360250255203225
$ws.Range("P9").Value = 0.001483950226849768
$ws.Range("Q9").Value = 0.2856166363535556
$ws.Range("R9").Value = 2.570549727182
$ws.Range("S9").Value = 0.0005764295076732777
$ws.Range("T9").Value = 0.0007277270118840214
$ws.Range("G10").Value = 14.273968
$ws.Range("H10").Value = 42.821904
$ws.Range("I10").Value = 0.09215331457489899
$ws.Range("J10").Value = 0.1066430931325968
$ws.Range("M10").Value = 0.2253945
$ws.Range("N10").Value = 0.450789
$ws.Range("O10").Value = 0.07045953565491116
$ws.Range("P10").Value = 0.05124470957648078
$ws.Range("Q10").Value = 3.217273880376
$ws.Range("R10").Value = 19.303643282256
$ws.Range("S10").Value = 0.006493079754008339
$ws.Range("T10").Value = 0.005464894335917514
$ws.Range("G11").Value = 14.273968
$ws.Range("H11").Value = 42.821904
$ws.Range("I11").Value = 0.09215331457489899
$ws.Range("J11").Value = 0.1066430931325968
$ws.Range("O11").Value = 0.7485640341142095
$ws.Range("P11").Value = 0.8166377943957063
$ws.Range("Q11").Value = 34.18040570888534
$ws.Range("R11").Value = 307.6236513799681
$ws.Range("S11").Value = 0.06898265691518217
$ws.Range("T11").Value = 0.08708878036333972
$ws.Range("G12").Value = 14.273968
$ws.Range("H12").Value = 42.821904
$ws.Range("I12").Value = 0.09215331457489899
$ws.Range("J12").Value = 0.1066430931325968
$ws.Range("M12").Value = 0.574578
$ws.Range("N12").Value = 1.149156
$ws.Range("O12").Value = 0.1796161799756762
$ws.Range("P12").Value = 0.1306335458009631
$ws.Range("Q12").Value = 8.201507985504001
$ws.Range("R12").Value = 49.20904791302401
$ws.Range("S12").Value = 0.01655222633604016
$ws.Range("T12").Value = 0.01393116539109345
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 14.273968
$ws.Range("H13").Value = 42.821904
$ws.Range("I13").Value = 0.09215331457489899
$ws.Range("J13").Value = 0.1066430931325968
$ws.Range("M13").Value = 0.004351333333333333
$ws.Range("N13").Value = 0.013054
$ws.Range("O13").Value = 0.001360250255203225
$ws.Range("P13").Value = 0.001483950226849768
$ws.Range("Q13").Value = 0.06211079275733333
$ws.Range("R13").Value = 0.558997134816
$ws.Range("S13").Value = 0.0001253515696683295
$ws.Range("T13").Value = 0.0001582530422460779
$ws.Range("G14").Value = 28.1345445
$ws.Range("H14").Value = 56.26908899999999
$ws.Range("I14").Value = 0.1816377569103415
$ws.Range("J14").Value = 0.1401317815927423
$ws.Range("M14").Value = 0.2253945
$ws.Range("N14").Value = 0.450789
$ws.Range("O14").Value = 0.07045953565491116
$ws.Range("P14").Value = 0.05124470957648078
$ws.Range("Q14").Value = 6.34137159030525
$ws.Range("R14").Value = 25.365486361221
$ws.Range("S14").Value = 0.01279811200930229
$ws.Range("T14").Value = 0.007181012450154912
$ws.Range("G15").Value = 28.1345445
$ws.Range("H15").Value = 56.26908899999999
$ws.Range("I15").Value = 0.1816377569103415
$ws.Range("J15").Value = 0.1401317815927423
$ws.Range("O15").Value = 0.7485640341142095
$ws.Range("P15").Value = 0.8166377943957063
$ws.Range("Q15").Value = 67.370905234248
$ws.Range("R15").Value = 404.225431405488
$ws.Range("S15").Value = 0.1359674920602614
$ws.Range("T15").Value = 0.1144369090446379
$ws.Range("G16").Value = 28.1345445
$ws.Range("H16").Value = 56.26908899999999
$ws.Range("I16").Value = 0.1816377569103415
$ws.Range("J16").Value = 0.1401317815927423
$ws.Range("M16").Value = 0.574578
$ws.Range("N16").Value = 1.149156
$ws.Range("O16").Value = 0.1796161799756762
$ws.Range("P16").Value = 0.1306335458009631
$ws.Range("Q16").Value = 16.165490309721
$ws.Range("R16").Value = 64.66196123888399
$ws.Range("S16").Value = 0.03262508003558602
$ws.Range("T16").Value = 0.01830591150886605
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("G17").Value = 28.1345445
$ws.Range("H17").Value = 56.26908899999999
$ws.Range("I17").Value = 0.1816377569103415
$ws.Range("J17").Value = 0.1401317815927423
$ws.Range("M17").Value = 0.004351333333333333
$ws.Range("N17").Value = 0.013054
$ws.Range("O17").Value = 0.001360250255203225
$ws.Range("P17").Value = 0.001483950226849768
$ws.Range("Q17").Value = 0.122422781301
$ws.Range("R17").Value = 0.7345366878059999
$ws.Range("S17").Value = 0.0002470728051918334
$ws.Range("T17").Value = 0.000207948589083412
$ws.Range("G18").Value = 2.253741
$ws.Range("H18").Value = 6.761222999999999
$ws.Range("I18").Value = 0.0145502430258599
$ws.Range("J18").Value = 0.01683805872058504
$ws.Range("M18").Value = 0.2253945
$ws.Range("N18").Value = 0.450789
$ws.Range("O18").Value = 0.07045953565491116
$ws.Range("P18").Value = 0.05124470957648078
$ws.Range("Q18").Value = 0.5079808258245
$ws.Range("R18").Value = 3.047884954947
$ws.Range("S18").Value = 0.001025203367268198
$ws.Range("T18").Value = 0.0008628614289681097
$ws.Range("G19").Value = 2.253741
$ws.Range("H19").Value = 6.761222999999999
$ws.Range("I19").Value = 0.0145502430258599
$ws.Range("J19").Value = 0.01683805872058504
$ws.Range("O19").Value = 0.7485640341142095
$ws.Range("P19").Value = 0.8166377943957063
$ws.Range("Q19").Value = 5.396802188624
$ws.Range("R19").Value = 48.57121969761599
$ws.Range("S19").Value = 0.01089178861677983
$ws.Range("T19").Value = 0.01375059513548395
$ws.Range("G20").Value = 2.253741
$ws.Range("H20").Value = 6.761222999999999
$ws.Range("I20").Value = 0.0145502430258599
$ws.Range("J20").Value = 0.01683805872058504
$ws.Range("M20").Value = 0.574578
$ws.Range("N20").Value = 1.149156
$ws.Range("O20").Value = 0.1796161799756762
$ws.Range("P20").Value = 0.1306335458009631
$ws.Range("Q20").Value = 1.294949996298
$ws.Range("R20").Value = 7.769699977788
$ws.Range("S20").Value = 0.002613459070022679
$ws.Range("T20").Value = 0.002199615315074851
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("G21").Value = 2.253741
$ws.Range("H21").Value = 6.761222999999999
$ws.Range("I21").Value = 0.0145502430258599
$ws.Range("J21").Value = 0.01683805872058504
$ws.Range("M21").Value = 0.004351333333333333
$ws.Range("N21").Value = 0.013054
$ws.Range("O21").Value = 0.001360250255203225
$ws.Range("P21").Value = 0.001483950226849768
$ws.Range("Q21").Value = 0.009806778337999998
$ws.Range("R21").Value = 0.08826100504199999
$ws.Range("S21").Value = 0.00001979197178919488
$ws.Range("T21").Value = 0.00002498684105812188
$ws.Range("G22").Value = 9.590102333333334
$ws.Range("H22").Value = 28.770307
$ws.Range("I22").Value = 0.06191408843911794
$ws.Range("J22").Value = 0.07164918516594687
$ws.Range("M22").Value = 0.2253945
$ws.Range("N22").Value = 0.450789
$ws.Range("O22").Value = 0.07045953565491116
$ws.Range("P22").Value = 0.05124470957648078
$ws.Range("Q22").Value = 2.1615563203705
$ws.Range("R22").Value = 12.969337922223
$ws.Range("S22").Value = 0.004362437921917353
$ws.Range("T22").Value = 0.003671641685220442
$ws.Range("G23").Value = 9.590102333333334
$ws.Range("H23").Value = 28.770307
$ws.Range("I23").Value = 0.06191408843911794
$ws.Range("J23").Value = 0.07164918516594687
$ws.Range("O23").Value = 0.7485640341142095
$ws.Range("P23").Value = 0.8166377943957063
$ws.Range("Q23").Value = 22.96443347379378
$ws.Range("R23").Value = 206.679901264144
$ws.Range("S23").Value = 0.04634665981049006
$ws.Range("T23").Value = 0.05851143254416841
$ws.Range("G24").Value = 9.590102333333334
$ws.Range("H24").Value = 28.770307
$ws.Range("I24").Value = 0.06191408843911794
$ws.Range("J24").Value = 0.07164918516594687
$ws.Range("M24").Value = 0.574578
$ws.Range("N24").Value = 1.149156
$ws.Range("O24").Value = 0.1796161799756762
$ws.Range("P24").Value = 0.1306335458009631
$ws.Range("Q24").Value = 5.510261818482001
$ws.Range("R24").Value = 33.061570910892
$ws.Range("S24").Value = 0.01112077205211054
$ws.Range("T24").Value = 0.009359787111977406
$ws.Range("D25").Value = "Resolving-Mac"
$ws.Range("G25").Value = 9.590102333333334
$ws.Range("H25").Value = 28.770307
$ws.Range("I25").Value = 0.06191408843911794
$ws.Range("J25").Value = 0.07164918516594687
$ws.Range("M25").Value = 0.004351333333333333
$ws.Range("N25").Value = 0.013054
$ws.Range("O25").Value = 0.001360250255203225
$ws.Range("P25").Value = 0.001483950226849768
$ws.Range("Q25").Value = 0.04172973195311111
$ws.Range("R25").Value = 0.375567587578
$ws.Range("S25").Value = 0.00008421865459998525
$ws.Range("T25").Value = 0.0001063238245806079
